# The commit updates the "Förändrad" (Changed) date column C for every
# data row (rows 2-292) from 2023-09-03 (serial 45172) to 2023-09-06
# (serial 45175). All other cells/styles/formulas are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 292; $r++) {
    $ws.Range("C$r").Value = 45175
}
